## "do not send empty reports"
## A container (FCIU8790621 / WASH / floor-board dirt damage, vessel
## "TR ARAMIS V 016N", rotation 2018/3965) that was previously being
## skipped because the shipment had no damage/remark rows is now
## included, so it gets appended as the last data row on both the
## "In Report" and "Stock Report" sheets (and their totals / summaries
## are bumped accordingly).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "In Report": append row 23 (14th container)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("In Report")

$ws.Range("A4").Value = "Total number of containers : 14"

$ws.Range("A23").Value = 14
$ws.Range("B23").Value = 845
$ws.Range("C23").Value = "FCIU8790621"
$ws.Range("D23").Value = 40
$ws.Range("E23").Value = "HC"
$ws.Range("F23").Value = "ISATL"
$ws.Range("G23").Value = "ISATL"
$ws.Range("H23").Value = "APL"
$ws.Range("I23").Value = "APL"
$ws.Range("J23").Value = "CPA"
$ws.Range("K23").Value = "TR ARAMIS V 016N"
$ws.Range("L23").Value = "2018/3965"
$ws.Range("M23").NumberFormat = "YYYY-MM-DD"
$ws.Range("M23").Value = Get-Date -Year 2018 -Month 12 -Day 26 -Hour 12 -Minute 44 -Second 20
$ws.Range("N23").Value = "WASH"
$ws.Range("O23").Value = ""
$ws.Range("P23").Value = ""
$ws.Range("Q23").NumberFormat = "YYYY-MM-DD"
$ws.Range("R23").Value = ""
$ws.Range("S23").Value = "FLOORS-()"
$ws.Range("T23").Value = "Floor board-()"
$ws.Range("U23").Value = "flooor board dirty by oil stain ,black stain ,& ink dirty "
$ws.Range("V23").Value = ""
$ws.Range("W23").Value = ""
$ws.Range("X23").Value = ""

# column U ("Damage Description") grows wide enough to hold the new
# remark text
$ws.Columns.Item(21).ColumnWidth = 39.7142857142857

# ---------------------------------------------------------------
# Sheet "In Report Summary": bump the 40-HC / APL / APL totals row
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("In Report Summary")
$ws.Range("H7").Value = 8
$ws.Range("K7").Value = 11
$ws.Range("L7").Value = 22

# ---------------------------------------------------------------
# Sheet "Stock Report": append row 28 (18th container)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Stock Report")

$ws.Range("A4").Value = "Total number of containers : 18"

$ws.Range("A28").Value = 18
$ws.Range("B28").Value = 845
$ws.Range("C28").Value = "FCIU8790621"
$ws.Range("D28").Value = 40
$ws.Range("E28").Value = "HC"
$ws.Range("F28").Value = "APL"
$ws.Range("G28").Value = "APL"
$ws.Range("H28").Value = "ISATL"
$ws.Range("I28").Value = "ISATL"
$ws.Range("J28").Value = "TR ARAMIS V 016N"
$ws.Range("K28").Value = "2018/3965"
$ws.Range("L28").Value = "CPA"
$ws.Range("M28").NumberFormat = "YYYY-MM-DD"
$ws.Range("M28").Value = Get-Date -Year 2018 -Month 12 -Day 26 -Hour 12 -Minute 44 -Second 20
$ws.Range("N28").Value = "WASH"
$ws.Range("O28").Value = ""
$ws.Range("P28").Value = 1
$ws.Range("Q28").Value = 0
$ws.Range("R28").Value = "Empty"
$ws.Range("S28").Value = ""
$ws.Range("T28").Value = ""
$ws.Range("U28").NumberFormat = "YYYY-MM-DD"
$ws.Range("V28").Value = ""
$ws.Range("W28").Value = "FLOORS-()"
$ws.Range("X28").Value = "Floor board-()"
$ws.Range("Y28").Value = "flooor board dirty by oil stain ,black stain ,& ink dirty "
$ws.Range("Z28").Value = ""
$ws.Range("AA28").Value = ""
$ws.Range("AB28").Value = ""

# column Y ("Damage Description") grows wide enough to hold the new
# remark text
$ws.Columns.Item(25).ColumnWidth = 39.7142857142857

# ---------------------------------------------------------------
# Sheet "Stock Report Summary": bump the 40-HC / APL / APL totals row
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Stock Report Summary")
$ws.Range("H7").Value = 8
$ws.Range("K7").Value = 15
$ws.Range("L7").Value = 30
